{"js": "// Update the two-digit-divided-by-one-digit practice problems to a new\n// set of numbers, matching each old \"XX\u00f7Y=\" expression to its replacement.\nconst replacements = [\n  [\"16\u00f76=\", \"40\u00f76=\"],\n  [\"79\u00f75=\", \"64\u00f73=\"],\n  [\"96\u00f74=\", \"52\u00f72=\"],\n  [\"72\u00f75=\", \"67\u00f73=\"],\n  [\"38\u00f76=\", \"84\u00f72=\"],\n  [\"33\u00f73=\", \"47\u00f72=\"],\n  [\"67\u00f78=\", \"75\u00f78=\"],\n  [\"76\u00f76=\", \"90\u00f77=\"],\n  [\"14\u00f73=\", \"86\u00f76=\"],\n  [\"32\u00f72=\", \"25\u00f76=\"],\n  [\"15\u00f79=\", \"12\u00f75=\"],\n  [\"81\u00f72=\", \"74\u00f77=\"],\n  [\"31\u00f72=\", \"64\u00f74=\"],\n  [\"50\u00f74=\", \"84\u00f74=\"],\n  [\"55\u00f78=\", \"43\u00f72=\"],\n  [\"12\u00f76=\", \"31\u00f79=\"],\n  [\"84\u00f75=\", \"13\u00f73=\"],\n  [\"83\u00f75=\", \"91\u00f75=\"],\n  [\"52\u00f77=\", \"65\u00f72=\"],\n  [\"68\u00f78=\", \"12\u00f73=\"],\n  [\"61\u00f72=\", \"20\u00f72=\"],\n  [\"29\u00f79=\", \"69\u00f72=\"],\n  [\"17\u00f72=\", \"68\u00f75=\"],\n  [\"38\u00f79=\", \"33\u00f77=\"],\n  [\"76\u00f78=\", \"69\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"16\u00f76=\", \"40\u00f76=\"),\n    @(\"79\u00f75=\", \"64\u00f73=\"),\n    @(\"96\u00f74=\", \"52\u00f72=\"),\n    @(\"72\u00f75=\", \"67\u00f73=\"),\n    @(\"38\u00f76=\", \"84\u00f72=\"),\n    @(\"33\u00f73=\", \"47\u00f72=\"),\n    @(\"67\u00f78=\", \"75\u00f78=\"),\n    @(\"76\u00f76=\", \"90\u00f77=\"),\n    @(\"14\u00f73=\", \"86\u00f76=\"),\n    @(\"32\u00f72=\", \"25\u00f76=\"),\n    @(\"15\u00f79=\", \"12\u00f75=\"),\n    @(\"81\u00f72=\", \"74\u00f77=\"),\n    @(\"31\u00f72=\", \"64\u00f74=\"),\n    @(\"50\u00f74=\", \"84\u00f74=\"),\n    @(\"55\u00f78=\", \"43\u00f72=\"),\n    @(\"12\u00f76=\", \"31\u00f79=\"),\n    @(\"84\u00f75=\", \"13\u00f73=\"),\n    @(\"83\u00f75=\", \"91\u00f75=\"),\n    @(\"52\u00f77=\", \"65\u00f72=\"),\n    @(\"68\u00f78=\", \"12\u00f73=\"),\n    @(\"61\u00f72=\", \"20\u00f72=\"),\n    @(\"29\u00f79=\", \"69\u00f72=\"),\n    @(\"17\u00f72=\", \"68\u00f75=\"),\n    @(\"38\u00f79=\", \"33\u00f77=\"),\n    @(\"76\u00f78=\", \"69\u00f73=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
